# Applies the "Add data for 2021-10-29" update to the carjacking-by-
# neighborhood-by-month workbook:
#   - rename the sheet / header / shared-string label from
#     "...through Oct 20" to "...through Oct 21"
#   - bump a handful of per-neighborhood monthly counts (including some
#     brand-new, previously-empty cells) to reflect the newly ingested day

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet + update the "through October N" label (sheet tab + B1 header cell)
$ws.Name = "Through 2021-10-21"
$ws.Range("B1").Value = "October 2021 (through October 21)"

# Row 2 - Garfield Park
$ws.Range("L2").Value = 15
$ws.Range("AP2").Value = 2

# Row 4 - North Lawndale
$ws.Range("AF4").Value = 2
$ws.Range("AP4").Value = 4

# Row 9 - Grand Crossing
$ws.Range("L9").Value = 1
$ws.Range("V9").Value = 3

# Row 10 - Roseland
$ws.Range("L10").Value = 2

# Row 12 - South Shore
$ws.Range("AP12").Value = 1

# Row 14 - West Town
$ws.Range("AP14").Value = 2
$ws.Range("AZ14").Value = 1

# Row 17 - Chicago Lawn
$ws.Range("AP17").Value = 1

# Row 18 - Little Village
$ws.Range("L18").Value = 3

# Row 21 - Lower West Side
$ws.Range("B21").Value = 3

# Row 22 - Grand Boulevard
$ws.Range("B22").Value = 4

# Row 32 - Edgewater
$ws.Range("AF32").Value = 2

# Row 42 - Washington Park
$ws.Range("V42").Value = 1

# Row 66 - Calumet Heights
$ws.Range("L66").Value = 2

# Row 82 - Mckinley Park
$ws.Range("L82").Value = 1

# Row 93 - Portage Park
$ws.Range("V93").Value = 2

# Row 97 - South Deering
$ws.Range("L97").Value = 2
